$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) First paragraph: remove the trailing "." run and the _GoBack bookmark
#    that sat there ("Negocios 2011." -> "Negocios 2011").
#    The _GoBack bookmark is relocated later (step 9) to the end of the
#    paragraph that ends with " @COPAIS" inside the BEGIN/SELECT block.
# ---------------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

$r = $d.Content
$found = $r.Find.Execute("Negocios 2011.", $true, $false, $false, $false, $false, $true, 1, $false, "Negocios 2011", 2)

# ---------------------------------------------------------------------------
# 2) Merge "Realizar los siguientes procesos " + "sql" + ". " into one run,
#    dropping the spellStart/spellEnd proofErr markers around "sql".
# ---------------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute("Realizar los siguientes procesos sql. ", $true, $false, $false, $false, $false, $true, 1, $false, "Realizar los siguientes procesos sql. ", 2)

# ---------------------------------------------------------------------------
# 3) Drop the gramStart/gramEnd proofErr markers around "CHAR(" in the
#    DECLARE @COPAIS CHAR(3) line.
# ---------------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute("CHAR(", $true, $false, $false, $false, $false, $true, 1, $false, "CHAR(", 2)

# ---------------------------------------------------------------------------
# 4) Merge the three runs "'" + "esp" + "'" into a single run "'esp'",
#    dropping the spellStart/spellEnd markers around "esp". A straight
#    find/replace with an apostrophe gets smart-quoted, and reassigning the
#    exact same text is a no-op for the run-merge, so first stamp a filler
#    value, then overwrite it with the real text via direct Range.Text
#    assignment (bypasses the autocorrect “smart quote” substitution).
# ---------------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute("'esp'")
$r.Text = "zzApostropheFillerzz"
$r2 = $d.Content
$found2 = $r2.Find.Execute("zzApostropheFillerzz")
$r2.Text = "'esp'"

# ---------------------------------------------------------------------------
# 5) Drop the gramStart/gramEnd proofErr markers around "EXISTS(".
# ---------------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute("EXISTS(", $true, $false, $false, $false, $false, $true, 1, $false, "EXISTS(", 2)

# ---------------------------------------------------------------------------
# 6) First "VENTAS.clientes" occurrence (inside IF EXISTS(...)):
#    - merge " " + "VENTAS" -> " VENTAS" (drop spellStart)
#    - merge "clientes" + "  " -> "clientes  " (drop spellEnd)
# ---------------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute("IF EXISTS(SELECT * FROM VENTAS.clientes  WHERE", $true, $false, $false, $false, $false, $true, 1, $false, "IF EXISTS(SELECT * FROM VENTAS.clientes  WHERE", 2)

# ---------------------------------------------------------------------------
# 7) Second "VENTAS.clientes" occurrence (inside BEGIN ... SELECT * FROM ...):
#    - merge " " + "VENTAS" -> " VENTAS" (drop spellStart + gramStart)
#    - merge "clientes" + "  " -> "clientes  " (drop spellEnd)
#    - drop the gramEnd marker that trailed "WHERE"
# ---------------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute("SELECT * FROM VENTAS.clientes  WHERE IDPAIS = @COPAIS", $true, $false, $false, $false, $false, $true, 1, $false, "SELECT * FROM VENTAS.clientes  WHERE IDPAIS = @COPAIS", 2)

# ---------------------------------------------------------------------------
# 8) Re-add the _GoBack bookmark, now collapsed at the end of the paragraph
#    that ends in " @COPAIS" (the SELECT ... WHERE IDPAIS = @COPAIS line
#    inside the BEGIN block).
# ---------------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute("SELECT * FROM VENTAS.clientes  WHERE IDPAIS = @COPAIS")
$target = $r.Duplicate
$target.Collapse(0)
$d.Bookmarks.Add("_GoBack", $target)
